$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1449.2609
$ws.Range("J17").Value = 1449.2609
$ws.Range("L17").Value = 4347.7827
$ws.Range("N17").Value = -4683.7827
$ws.Range("H19").Value = 1602.7391
$ws.Range("I19").Value = 1458.5555
$ws.Range("J19").Value = 1695.4286
$ws.Range("K19").Value = 1458.5555
$ws.Range("L19").Value = 1695.4286
$ws.Range("M19").Value = -1283.5555
$ws.Range("N19").Value = -2045.4286
$ws.Range("H87").Value = 27000
$ws.Range("J87").Value = 27000
$ws.Range("L87").Value = 27000
$ws.Range("N87").Value = -29496
$ws.Range("H90").Value = 27000
$ws.Range("J90").Value = 27000
$ws.Range("L90").Value = 81000
$ws.Range("N90").Value = -93480
$ws.Range("H121").Value = 1162
$ws.Range("J121").Value = 1162
$ws.Range("L121").Value = 3486
$ws.Range("N121").Value = -6980
$ws.Range("H132").Value = 1118.5834
$ws.Range("I132").Value = 1102.1364
$ws.Range("J132").Value = 1299.5
$ws.Range("K132").Value = 3306.4092
$ws.Range("L132").Value = 3898.5
$ws.Range("M132").Value = -776.4092
$ws.Range("N132").Value = -8958.5
$ws.Range("H137").Value = 39073.37
$ws.Range("I137").Value = 1037.6666
$ws.Range("J137").Value = 43827.832
$ws.Range("K137").Value = 3112.9998
$ws.Range("L137").Value = 131483.496
$ws.Range("M137").Value = -562.9998
$ws.Range("N137").Value = -136583.496
$ws.Range("H138").Value = 2602.6924
$ws.Range("I138").Value = 2802.9333
$ws.Range("J138").Value = 2477.5417
$ws.Range("K138").Value = 8408.7999
$ws.Range("L138").Value = 7432.625100000001
$ws.Range("M138").Value = -3268.7999
$ws.Range("N138").Value = -17712.6251
$ws.Range("H141").Value = 1169106.4
$ws.Range("I141").Value = 1649032.6
$ws.Range("J141").Value = 3571.1428
$ws.Range("K141").Value = 4947097.800000001
$ws.Range("L141").Value = 10713.4284
$ws.Range("M141").Value = -4941917.800000001
$ws.Range("N141").Value = -21073.4284

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2308.4707
$ws.Range("I32").Value = 1800.3148
$ws.Range("J32").Value = 4268.5
$ws.Range("K32").Value = 1800.3148
$ws.Range("L32").Value = 4268.5
$ws.Range("M32").Value = -1513.3148
$ws.Range("N32").Value = -4842.5
$ws.Range("H45").Value = 3090.1428
$ws.Range("I45").Value = 4324.25
$ws.Range("K45").Value = 4324.25
$ws.Range("M45").Value = -3947.25
$ws.Range("H61").Value = 3332.5557
$ws.Range("I61").Value = 2579.6667
$ws.Range("J61").Value = 4273.6665
$ws.Range("K61").Value = 2579.6667
$ws.Range("L61").Value = 4273.6665
$ws.Range("M61").Value = -2367.6667
$ws.Range("N61").Value = -4697.6665
$ws.Range("H74").Value = 2151
$ws.Range("I74").Value = 742
$ws.Range("K74").Value = 742
$ws.Range("M74").Value = 132
$ws.Range("H77").Value = 2151
$ws.Range("I77").Value = 742
$ws.Range("K77").Value = 3710
$ws.Range("M77").Value = 658
$ws.Range("H97").Value = 2138
$ws.Range("I97").Value = 1965.6
$ws.Range("J97").Value = 3000
$ws.Range("K97").Value = 1965.6
$ws.Range("L97").Value = 3000
$ws.Range("M97").Value = -1469.6
$ws.Range("N97").Value = -3992
$ws.Range("H122").Value = 33040.25
$ws.Range("I122").Value = 64728.168
$ws.Range("J122").Value = 1352.3334
$ws.Range("K122").Value = 194184.504
$ws.Range("L122").Value = 4057.0002
$ws.Range("M122").Value = -191734.504
$ws.Range("N122").Value = -8957.0002
$ws.Range("H132").Value = 2068.8572
$ws.Range("I132").Value = 1846.7709
$ws.Range("K132").Value = 5540.3127
$ws.Range("M132").Value = -3010.3127
$ws.Range("H136").Value = 3332.5557
$ws.Range("I136").Value = 2579.6667
$ws.Range("J136").Value = 4273.6665
$ws.Range("K136").Value = 7739.000100000001
$ws.Range("L136").Value = 12820.9995
$ws.Range("M136").Value = -5189.000100000001
$ws.Range("N136").Value = -17920.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1959.9
$ws.Range("I20").Value = 1920.5333
$ws.Range("J20").Value = 2078
$ws.Range("K20").Value = 1920.5333
$ws.Range("L20").Value = 2078
$ws.Range("M20").Value = -1673.5333
$ws.Range("N20").Value = -2572
$ws.Range("H86").Value = 601701.2
$ws.Range("I86").Value = 2835.3333
$ws.Range("J86").Value = 1500000
$ws.Range("K86").Value = 2835.3333
$ws.Range("L86").Value = 1500000
$ws.Range("M86").Value = -1712.3333
$ws.Range("N86").Value = -1502246
$ws.Range("H88").Value = 34432.332
$ws.Range("J88").Value = 34432.332
$ws.Range("L88").Value = 34432.332
$ws.Range("N88").Value = -35244.332
$ws.Range("H89").Value = 601701.2
$ws.Range("I89").Value = 2835.3333
$ws.Range("J89").Value = 1500000
$ws.Range("K89").Value = 14176.6665
$ws.Range("L89").Value = 7500000
$ws.Range("M89").Value = -8560.6665
$ws.Range("N89").Value = -7511232
$ws.Range("H91").Value = 34432.332
$ws.Range("J91").Value = 34432.332
$ws.Range("L91").Value = 34432.332
$ws.Range("N91").Value = -37240.332
$ws.Range("H94").Value = 1180.4
$ws.Range("I94").Value = 971.6842
$ws.Range("J94").Value = 1540.909
$ws.Range("K94").Value = 971.6842
$ws.Range("L94").Value = 1540.909
$ws.Range("M94").Value = -520.6842
$ws.Range("N94").Value = -2442.909
$ws.Range("H95").Value = 67945
$ws.Range("J95").Value = 67945
$ws.Range("L95").Value = 67945
$ws.Range("N95").Value = -73437
$ws.Range("H134").Value = 12741.1
$ws.Range("I134").Value = 14938.875
$ws.Range("K134").Value = 44816.625
$ws.Range("M134").Value = -42281.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 421.25
$ws.Range("I10").Value = 421.25
$ws.Range("K10").Value = 421.25
$ws.Range("M10").Value = -282.25
$ws.Range("H16").Value = 1105.25
$ws.Range("I16").Value = 1105.25
$ws.Range("K16").Value = 1105.25
$ws.Range("M16").Value = -818.25
$ws.Range("H31").Value = 1291.3414
$ws.Range("I31").Value = 830.5417
$ws.Range("J31").Value = 1941.8823
$ws.Range("K31").Value = 830.5417
$ws.Range("L31").Value = 1941.8823
$ws.Range("M31").Value = -535.5417
$ws.Range("N31").Value = -2531.8823
$ws.Range("H34").Value = 1291.3414
$ws.Range("I34").Value = 830.5417
$ws.Range("J34").Value = 1941.8823
$ws.Range("K34").Value = 830.5417
$ws.Range("L34").Value = 1941.8823
$ws.Range("M34").Value = -628.5417
$ws.Range("N34").Value = -2345.8823
$ws.Range("H58").Value = 2899748.5
$ws.Range("I58").Value = 5435847
$ws.Range("J58").Value = 1350.4286
$ws.Range("K58").Value = 5435847
$ws.Range("L58").Value = 1350.4286
$ws.Range("M58").Value = -5435644
$ws.Range("N58").Value = -1756.4286
$ws.Range("H74").Value = 29999.666
$ws.Range("J74").Value = 29999.666
$ws.Range("L74").Value = 29999.666
$ws.Range("N74").Value = -31747.666
$ws.Range("H77").Value = 29999.666
$ws.Range("J77").Value = 29999.666
$ws.Range("L77").Value = 89998.998
$ws.Range("N77").Value = -98734.998
$ws.Range("H113").Value = 1105.25
$ws.Range("I113").Value = 1105.25
$ws.Range("K113").Value = 1105.25
$ws.Range("M113").Value = 1064.75
$ws.Range("H132").Value = 2683.7407
$ws.Range("I132").Value = 2250.647
$ws.Range("K132").Value = 6751.941
$ws.Range("M132").Value = -4221.941
$ws.Range("H134").Value = 1538.5625
$ws.Range("I134").Value = 829
$ws.Range("K134").Value = 2487
$ws.Range("M134").Value = 48
$ws.Range("H136").Value = 2899748.5
$ws.Range("I136").Value = 5435847
$ws.Range("J136").Value = 1350.4286
$ws.Range("K136").Value = 16307541
$ws.Range("L136").Value = 4051.2858
$ws.Range("M136").Value = -16304991
$ws.Range("N136").Value = -9151.2858

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5333.3335
$ws.Range("I70").Value = 7000
$ws.Range("J70").Value = 4500
$ws.Range("K70").Value = 7000
$ws.Range("L70").Value = 4500
$ws.Range("M70").Value = -6730
$ws.Range("N70").Value = -5040
$ws.Range("H73").Value = 5333.3335
$ws.Range("I73").Value = 7000
$ws.Range("J73").Value = 4500
$ws.Range("K73").Value = 7000
$ws.Range("L73").Value = 4500
$ws.Range("M73").Value = -6064
$ws.Range("N73").Value = -6372
$ws.Range("H126").Value = 2418634
$ws.Range("I126").Value = 3706984.5
$ws.Range("K126").Value = 11120953.5
$ws.Range("M126").Value = -11118483.5
$ws.Range("H132").Value = 1042068.7
$ws.Range("I132").Value = 1604275.2
$ws.Range("J132").Value = 4148.923
$ws.Range("K132").Value = 4812825.6
$ws.Range("L132").Value = 12446.769
$ws.Range("M132").Value = -4810295.6
$ws.Range("N132").Value = -17506.769

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3683
$ws.Range("I22").Value = 1325
$ws.Range("J22").Value = 5255
$ws.Range("K22").Value = 1325
$ws.Range("L22").Value = 5255
$ws.Range("M22").Value = -1030
$ws.Range("N22").Value = -5845
$ws.Range("H27").Value = 3683
$ws.Range("I27").Value = 1325
$ws.Range("J27").Value = 5255
$ws.Range("K27").Value = 1325
$ws.Range("L27").Value = 5255
$ws.Range("M27").Value = -1218
$ws.Range("N27").Value = -5469
$ws.Range("H122").Value = 11555.667
$ws.Range("I122").Value = 10667.333
$ws.Range("K122").Value = 32001.999
$ws.Range("M122").Value = -29551.999
$ws.Range("H132").Value = 3178.0625
$ws.Range("I132").Value = 973.53845
$ws.Range("J132").Value = 4686.421
$ws.Range("K132").Value = 2920.61535
$ws.Range("L132").Value = 14059.263
$ws.Range("M132").Value = -390.61535
$ws.Range("N132").Value = -19119.263

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3009.2632
$ws.Range("I132").Value = 2573.1667
$ws.Range("J132").Value = 3756.8572
$ws.Range("K132").Value = 7719.500100000001
$ws.Range("L132").Value = 11270.5716
$ws.Range("M132").Value = -5189.500100000001
$ws.Range("N132").Value = -16330.5716
$ws.Range("H136").Value = 25256134
$ws.Range("I136").Value = 50508268
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 151524804
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = -151522254
$ws.Range("N136").Value = -17100
